# Apply numeric updates scraped from the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H33").Value = 9616186
$ws_ALC.Range("I33").Value = 12500407
$ws_ALC.Range("K33").Value = 12500407
$ws_ALC.Range("M33").Value = -12500178
$ws_ALC.Range("H40").Value = 1245
$ws_ALC.Range("J40").Value = 1245
$ws_ALC.Range("L40").Value = 1245
$ws_ALC.Range("N40").Value = -1595
$ws_ALC.Range("H41").Value = 649.25
$ws_ALC.Range("J41").Value = 2000
$ws_ALC.Range("L41").Value = 2000
$ws_ALC.Range("N41").Value = -2880
$ws_ALC.Range("H43").Value = 8238.875
$ws_ALC.Range("I43").Value = 8000
$ws_ALC.Range("J43").Value = 8382.200000000001
$ws_ALC.Range("K43").Value = 8000
$ws_ALC.Range("L43").Value = 8382.200000000001
$ws_ALC.Range("M43").Value = -7931
$ws_ALC.Range("N43").Value = -8520.200000000001
$ws_ALC.Range("H106").Value = 2605.45
$ws_ALC.Range("J106").Value = 3666.3333
$ws_ALC.Range("L106").Value = 3666.3333
$ws_ALC.Range("N106").Value = -4928.3333
$ws_ALC.Range("H112").Value = 84852.69500000001
$ws_ALC.Range("J112").Value = 93766.82000000001
$ws_ALC.Range("L112").Value = 281300.46
$ws_ALC.Range("N112").Value = -283516.46
$ws_ALC.Range("H132").Value = 9721.317999999999
$ws_ALC.Range("I132").Value = 3124.6
$ws_ALC.Range("J132").Value = 23857.143
$ws_ALC.Range("K132").Value = 9373.799999999999
$ws_ALC.Range("L132").Value = 71571.429
$ws_ALC.Range("M132").Value = -6843.799999999999
$ws_ALC.Range("N132").Value = -76631.429
$ws_ALC.Range("H135").Value = 1197.8334
$ws_ALC.Range("J135").Value = 3000
$ws_ALC.Range("L135").Value = 27000
$ws_ALC.Range("N135").Value = -32070
$ws_ALC.Range("H138").Value = 4219.3335
$ws_ALC.Range("I138").Value = 8763.200000000001
$ws_ALC.Range("J138").Value = 3605.2974
$ws_ALC.Range("K138").Value = 26289.6
$ws_ALC.Range("L138").Value = 10815.8922
$ws_ALC.Range("M138").Value = -21149.6
$ws_ALC.Range("N138").Value = -21095.8922
$ws_ALC.Range("H139").Value = 92249.625
$ws_ALC.Range("J139").Value = 93999.57000000001
$ws_ALC.Range("L139").Value = 93999.57000000001
$ws_ALC.Range("N139").Value = -104279.57
$ws_ALC.Range("H141").Value = 3401.3572
$ws_ALC.Range("I141").Value = 2963.0908
$ws_ALC.Range("K141").Value = 8889.2724
$ws_ALC.Range("M141").Value = -3709.2724
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 3387.7354
$ws_ARM.Range("I32").Value = 1512.4615
$ws_ARM.Range("K32").Value = 1512.4615
$ws_ARM.Range("M32").Value = -1225.4615
$ws_ARM.Range("H63").Value = 3540.3333
$ws_ARM.Range("H66").Value = 3540.3333
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H19").Value = 12144.143
$ws_BSM.Range("I19").Value = 12144.143
$ws_BSM.Range("J19").Value = 0
$ws_BSM.Range("K19").Value = 12144.143
$ws_BSM.Range("L19").Value = 0
$ws_BSM.Range("N19").ClearContents()
$ws_BSM.Range("M19").Value = -11971.143
$ws_BSM.Range("H82").Value = 9825.833000000001
$ws_BSM.Range("I82").Value = 9825.833000000001
$ws_BSM.Range("J82").Value = 0
$ws_BSM.Range("K82").Value = 9825.833000000001
$ws_BSM.Range("L82").Value = 0
$ws_BSM.Range("M82").Value = -9442.833000000001
$ws_BSM.Range("N82").ClearContents()
$ws_BSM.Range("H85").Value = 9825.833000000001
$ws_BSM.Range("I85").Value = 9825.833000000001
$ws_BSM.Range("J85").Value = 0
$ws_BSM.Range("K85").Value = 9825.833000000001
$ws_BSM.Range("L85").Value = 0
$ws_BSM.Range("M85").Value = -8499.833000000001
$ws_BSM.Range("N85").ClearContents()
$ws_BSM.Range("H86").Value = 895.4286
$ws_BSM.Range("I86").Value = 763.8
$ws_BSM.Range("J86").Value = 1224.5
$ws_BSM.Range("K86").Value = 763.8
$ws_BSM.Range("L86").Value = 1224.5
$ws_BSM.Range("M86").Value = 359.2
$ws_BSM.Range("N86").Value = -3470.5
$ws_BSM.Range("H89").Value = 895.4286
$ws_BSM.Range("I89").Value = 763.8
$ws_BSM.Range("J89").Value = 1224.5
$ws_BSM.Range("K89").Value = 3819
$ws_BSM.Range("L89").Value = 6122.5
$ws_BSM.Range("M89").Value = 1797
$ws_BSM.Range("N89").Value = -17354.5
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H7").Value = 542.6667
$ws_CRP.Range("I7").Value = 57
$ws_CRP.Range("K7").Value = 57
$ws_CRP.Range("M7").Value = 56
$ws_CRP.Range("H31").Value = 111053.55
$ws_CRP.Range("I31").Value = 200965.1
$ws_CRP.Range("K31").Value = 200965.1
$ws_CRP.Range("M31").Value = -200670.1
$ws_CRP.Range("H34").Value = 111053.55
$ws_CRP.Range("I34").Value = 200965.1
$ws_CRP.Range("K34").Value = 200965.1
$ws_CRP.Range("M34").Value = -200763.1
$ws_CRP.Range("H132").Value = 64349764
$ws_CRP.Range("I132").Value = 2871
$ws_CRP.Range("K132").Value = 8613
$ws_CRP.Range("M132").Value = -6083
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H4").Value = 3323900.5
$ws_CUL.Range("I4").Value = 4525054.5
$ws_CUL.Range("K4").Value = 13575163.5
$ws_CUL.Range("M4").Value = -13575051.5
$ws_CUL.Range("H14").Value = 96.75
$ws_CUL.Range("I14").Value = 96.75
$ws_CUL.Range("K14").Value = 290.25
$ws_CUL.Range("M14").Value = -117.25
$ws_CUL.Range("H20").Value = 2000
$ws_CUL.Range("J20").Value = 5000
$ws_CUL.Range("L20").Value = 15000
$ws_CUL.Range("N20").Value = -15454
$ws_CUL.Range("H58").Value = 6666
$ws_CUL.Range("J58").Value = 7499
$ws_CUL.Range("L58").Value = 22497
$ws_CUL.Range("N58").Value = -22753
$ws_CUL.Range("H131").Value = 1573.15
$ws_CUL.Range("I131").Value = 1474
$ws_CUL.Range("J131").Value = 1575.1735
$ws_CUL.Range("K131").Value = 4422
$ws_CUL.Range("L131").Value = 4725.520500000001
$ws_CUL.Range("M131").Value = 618
$ws_CUL.Range("N131").Value = -14805.5205
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H70").Value = 9140.117
$ws_GSM.Range("I70").Value = 8718.130999999999
$ws_GSM.Range("J70").Value = 10022.454
$ws_GSM.Range("K70").Value = 8718.130999999999
$ws_GSM.Range("L70").Value = 10022.454
$ws_GSM.Range("M70").Value = -8448.130999999999
$ws_GSM.Range("N70").Value = -10562.454
$ws_GSM.Range("H73").Value = 9140.117
$ws_GSM.Range("I73").Value = 8718.130999999999
$ws_GSM.Range("J73").Value = 10022.454
$ws_GSM.Range("K73").Value = 8718.130999999999
$ws_GSM.Range("L73").Value = 10022.454
$ws_GSM.Range("M73").Value = -7782.130999999999
$ws_GSM.Range("N73").Value = -11894.454
$ws_GSM.Range("H80").Value = 2637.8
$ws_GSM.Range("I80").Value = 1619.75
$ws_GSM.Range("K80").Value = 1619.75
$ws_GSM.Range("M80").Value = -621.75
$ws_GSM.Range("H83").Value = 2637.8
$ws_GSM.Range("I83").Value = 1619.75
$ws_GSM.Range("K83").Value = 8098.75
$ws_GSM.Range("M83").Value = -3106.75
$ws_GSM.Range("H132").Value = 553531.3
$ws_GSM.Range("I132").Value = 3809.3635
$ws_GSM.Range("K132").Value = 11428.0905
$ws_GSM.Range("M132").Value = -8898.0905
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H16").Value = 1696.7241
$ws_LTW.Range("J16").Value = 1809.4445
$ws_LTW.Range("L16").Value = 1809.4445
$ws_LTW.Range("N16").Value = -2149.4445
$ws_LTW.Range("H22").Value = 604.5
$ws_LTW.Range("I22").Value = 376.4
$ws_LTW.Range("J22").Value = 1174.75
$ws_LTW.Range("K22").Value = 376.4
$ws_LTW.Range("L22").Value = 1174.75
$ws_LTW.Range("M22").Value = -81.39999999999998
$ws_LTW.Range("N22").Value = -1764.75
$ws_LTW.Range("H27").Value = 604.5
$ws_LTW.Range("I27").Value = 376.4
$ws_LTW.Range("J27").Value = 1174.75
$ws_LTW.Range("K27").Value = 376.4
$ws_LTW.Range("L27").Value = 1174.75
$ws_LTW.Range("M27").Value = -269.4
$ws_LTW.Range("N27").Value = -1388.75
$ws_LTW.Range("H46").Value = 3182.3684
$ws_LTW.Range("I46").Value = 998.3333
$ws_LTW.Range("K46").Value = 998.3333
$ws_LTW.Range("M46").Value = -810.3333
$ws_LTW.Range("H132").Value = 1200616.1
$ws_LTW.Range("I132").Value = 4803
$ws_LTW.Range("J132").Value = 2722560
$ws_LTW.Range("K132").Value = 14409
$ws_LTW.Range("L132").Value = 8167680
$ws_LTW.Range("M132").Value = -11879
$ws_LTW.Range("N132").Value = -8172740
$ws_LTW.Range("H136").Value = 1603613.9
$ws_LTW.Range("I136").Value = 69914
$ws_LTW.Range("K136").Value = 209742
$ws_LTW.Range("M136").Value = -207192
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H54").Value = 12069.375
$ws_WVR.Range("I54").Value = 2760.8333
$ws_WVR.Range("K54").Value = 2760.8333
$ws_WVR.Range("M54").Value = -2240.8333
$ws_WVR.Range("H81").Value = 2074.5833
$ws_WVR.Range("I81").Value = 1310.4445
$ws_WVR.Range("J81").Value = 4367
$ws_WVR.Range("K81").Value = 2620.889
$ws_WVR.Range("L81").Value = 8734
$ws_WVR.Range("M81").Value = -1559.889
$ws_WVR.Range("N81").Value = -10856
$ws_WVR.Range("H84").Value = 2074.5833
$ws_WVR.Range("I84").Value = 1310.4445
$ws_WVR.Range("J84").Value = 4367
$ws_WVR.Range("K84").Value = 13104.445
$ws_WVR.Range("L84").Value = 43670
$ws_WVR.Range("M84").Value = -7800.445
$ws_WVR.Range("N84").Value = -54278
$ws_WVR.Range("H109").Value = 22666.666
$ws_WVR.Range("J109").Value = 22666.666
$ws_WVR.Range("L109").Value = 22666.666
$ws_WVR.Range("N109").Value = -25440.666
$ws_WVR.Range("H136").Value = 486581.8
$ws_WVR.Range("I136").Value = 7283
$ws_WVR.Range("J136").Value = 710254.6
$ws_WVR.Range("K136").Value = 21849
$ws_WVR.Range("L136").Value = 2130763.8
$ws_WVR.Range("M136").Value = -19299
$ws_WVR.Range("N136").Value = -2135863.8

Write-Output "Applied 221 cell updates across 8 sheets."
